$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text could be misread as a number/date by Excel's
# automatic type inference get a temporary Text number format so the
# COM layer stores them as strings (matching the source inlineStr cells),
# then the format stamp is cleared so no stray style index is introduced.
function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Row 2
$ws.Range("D2").Value = '27.416.36'
$ws.Range("E2").Value = '  -1.31%  '

# Row 3
$ws.Range("D3").Value = '1.830.77'
$ws.Range("E3").Value = '  -1.45%  '

# Row 4
Set-TextValue "D4" '1.006'
$ws.Range("E4").Value = '  -2.65%  '

# Row 5
Set-TextValue "D5" '314.91'
$ws.Range("E5").Value = '  -2.43%  '

# Row 6
$ws.Range("E6").Value = '  -2.79%  '

# Row 7
Set-TextValue "D7" '0.4295'
$ws.Range("E7").Value = '  -2.46%  '

# Row 8
Set-TextValue "D8" '0.3703'
$ws.Range("E8").Value = '  -2.89%  '

# Row 9
Set-TextValue "D9" '0.07265'
$ws.Range("E9").Value = '  -2.06%  '

# Row 10
Set-TextValue "D10" '0.8666'
$ws.Range("E10").Value = '  -2.33%  '

# Row 11
Set-TextValue "D11" '21.19'
$ws.Range("E11").Value = '  -1.79%  '

# Row 12
$ws.Range("D12").Value = '1.819.97'
$ws.Range("E12").Value = '  -1.88%  '

# Row 13
Set-TextValue "D13" '6.685'
$ws.Range("E13").Value = '  -0.55%  '

# Row 14
Set-TextValue "D14" '5.362'
$ws.Range("E14").Value = '  -2.86%  '

# Row 15
Set-TextValue "D15" '0.07063'
$ws.Range("E15").Value = '  -1.80%  '

# Row 16
Set-TextValue "D16" '87.92'
$ws.Range("E16").Value = '  +3.11%  '

# Row 17
$ws.Range("E17").Value = '  -3.13%  '

# Row 18
Set-TextValue "D18" '0.000008906'
$ws.Range("E18").Value = '  -1.92%  '

# Row 19
$ws.Range("E19").Value = '  -2.81%  '

# Row 20
Set-TextValue "D20" '15.24'
$ws.Range("E20").Value = '  -1.90%  '

# Row 21
$ws.Range("D21").Value = '27.456.20'
$ws.Range("E21").Value = '  -1.25%  '

# Row 22
Set-TextValue "D22" '5.169'
$ws.Range("E22").Value = '  -2.19%  '

# Row 23
Set-TextValue "D23" '10.92'
$ws.Range("E23").Value = '  -3.16%  '

# Row 24
$ws.Range("D24").Value = '2.057.18'
$ws.Range("E24").Value = '  -1.53%  '

# Row 25
$ws.Range("E25").Value = '  -3.18%  '

# Row 26
Set-TextValue "D26" '153.52'
$ws.Range("E26").Value = '  -3.66%  '

# Row 27
Set-TextValue "D27" '18.51'
$ws.Range("E27").Value = '  -1.27%  '

# Row 28
Set-TextValue "D28" '2.157'
$ws.Range("E28").Value = '  +7.88%  '

# Row 29
Set-TextValue "D29" '5.285'
$ws.Range("E29").Value = '  -1.30%  '

# Row 30
Set-TextValue "D30" '117.35'
$ws.Range("E30").Value = '  -0.64%  '

# Row 31
Set-TextValue "D31" '0.08860'
$ws.Range("E31").Value = '  -2.53%  '

# Row 32
Set-TextValue "D32" '1.209'
$ws.Range("E32").Value = '  -0.47%  '

# Row 33
Set-TextValue "D33" '0.7678'
$ws.Range("E33").Value = '  -0.82%  '

# Row 34
$ws.Range("E34").Value = '  -2.15%  '

# Row 35
Set-TextValue "D35" '2.898'
$ws.Range("E35").Value = '  -3.71%  '

# Row 36
$ws.Range("E36").Value = '  -2.99%  '

# Row 37
Set-TextValue "D37" '1.122'
$ws.Range("E37").Value = '  -2.61%  '

# Row 38
Set-TextValue "D38" '0.01962'
$ws.Range("E38").Value = '  -0.95%  '

# Row 39
Set-TextValue "D39" '0.05285'
$ws.Range("E39").Value = '  -0.21%  '

# Row 40
Set-TextValue "D40" '7.177'
$ws.Range("E40").Value = '  +3.94%  '

# Row 41
Set-TextValue "D41" '2.868'
$ws.Range("E41").Value = '  +0.12%  '

# Row 42
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue "D42" '0.1679'
$ws.Range("E42").Value = '  +0.35%  '

# Row 43
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue "D43" '0.5092'
$ws.Range("E43").Value = '  -1.97%  '

# Row 44
Set-TextValue "D44" '8.679'
$ws.Range("E44").Value = '  -0.95%  '

# Row 45
$ws.Range("E45").Value = '  -1.89%  '

# Row 46
Set-TextValue "D46" '0.4738'
$ws.Range("E46").Value = '  +0.42%  '

# Row 47
$ws.Range("E47").Value = '  -4.07%  '

# Row 48
Set-TextValue "D48" '0.06424'
$ws.Range("E48").Value = '  -2.42%  '

# Row 49
$ws.Range("E49").Value = '  -3.17%  '

# Row 50
Set-TextValue "D50" '1.668'
$ws.Range("E50").Value = '  -2.59%  '

# Row 51
Set-TextValue "D51" '1.829'
$ws.Range("E51").Value = '  -3.73%  '
